# "Generate Report for handback" - refresh the handback-status report with
# newly generated handoff/handback timestamps for the 466af032... entry
# (row 4 of the per-language sheets) on both the zh-cn and de-de sheets.
#
# Column D = "Correspond Handoff Datetime"
# Column G = "Correspond Handback DateTime"

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D4").Value = "2016-01-19 05:38:11"
$wsZh.Range("G4").Value = "2016-01-19 05:38:56"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D4").Value = "2016-01-19 05:38:22"
$wsDe.Range("G4").Value = "2016-01-19 05:39:13"
